$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 151 (pushes existing rows 151:208 down to 152:209,
# and extends the used range to A1:R209).
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A151").Value = 6
$ws.Range("B151").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C151").Value = "Metropolitana"
$ws.Range("D151").Value = 44468
$ws.Range("E151").Value = 13
$ws.Range("F151").Value = 100112032
$ws.Range("G151").Value = "Zapallo italiano"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 150
$ws.Range("K151").Value = 18000
$ws.Range("L151").Value = 19000
$ws.Range("M151").Value = 18333
$ws.Range("N151").Value = "`$/caja 50 unidades"
$ws.Range("O151").Value = "Región de Arica y Parinacota"
$ws.Range("P151").Value = 367
$ws.Range("Q151").Value = 50
$ws.Range("R151").Value = "Hortaliza"
